# Updated cryptos list on Thu May  2 12:40:47 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddr, $text) {
    $r = $ws.Range($rangeAddr)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.ClearFormats()
}

# Row 2 - Bitcoin
Set-TextValue "D2" "58.840.49"
Set-TextValue "E2" "  +1.01%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.996.28"
Set-TextValue "E3" "  +2.21%  "

# Row 4 - TetherUSD
Set-TextValue "E4" "  -0.10%  "

# Row 5 - BNB
Set-TextValue "D5" "562.22"
Set-TextValue "E5" "  +1.24%  "

# Row 6 - Solana
Set-TextValue "D6" "137.61"
Set-TextValue "E6" "  +10.45%  "

# Row 7 - USDC
Set-TextValue "E7" "  +0.00%  "

# Row 8 - XRP
Set-TextValue "D8" "0.517"
Set-TextValue "E8" "  +4.25%  "

# Row 9 - LidoStakedEther
Set-TextValue "D9" "2.991.27"
Set-TextValue "E9" "  +2.20%  "

# Row 10 - Dogecoin
Set-TextValue "E10" "  +3.91%  "

# Row 11 - Toncoin
Set-TextValue "E11" "  +1.99%  "

# Row 12 - Cardano
Set-TextValue "E12" "  +2.84%  "

# Row 13 - ShibaInu
Set-TextValue "E13" "  +6.60%  "

# Row 14 - Avalanche
Set-TextValue "D14" "33.71"
Set-TextValue "E14" "  +4.39%  "

# Row 15 - TRON
Set-TextValue "E15" "  +2.60%  "

# Row 16 - WrappedliquidstakedEther2.0
Set-TextValue "D16" "3.490.94"
Set-TextValue "E16" "  +2.33%  "

# Row 17 - Polkadot
Set-TextValue "D17" "7.00"
Set-TextValue "E17" "  +6.45%  "

# Row 18 - WrappedEther
Set-TextValue "D18" "2.992.71"
Set-TextValue "E18" "  +2.27%  "

# Row 19 - WrappedBTC
Set-TextValue "D19" "58.849.49"
Set-TextValue "E19" "  +1.22%  "

# Row 20 - BitcoinCash
Set-TextValue "D20" "426.36"
Set-TextValue "E20" "  +3.22%  "

# Row 21 - Chainlink
Set-TextValue "D21" "13.61"
Set-TextValue "E21" "  +4.81%  "

# Row 22 - Polygon
Set-TextValue "D22" "0.713"
Set-TextValue "E22" "  +6.98%  "

# Row 23 - was InternetComputer(DFINITY), now Uniswap
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue "D23" "7.15"
Set-TextValue "E23" "  +3.61%  "

# Row 24 - was Uniswap, now InternetComputer(DFINITY)
$ws.Range("B24").Value = "InternetComputer(DFINITY)"
$ws.Range("C24").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D24" "13.46"
Set-TextValue "E24" "  +4.17%  "

# Row 25 - Litecoin
Set-TextValue "D25" "80.34"
Set-TextValue "E25" "  +3.66%  "

# Row 26 - Dai
Set-TextValue "E26" "  -0.13%  "

# Row 27 - FirstDigitalUSD
Set-TextValue "D27" "0.999"
Set-TextValue "E27" "  -0.04%  "

# Row 28 - ImmutableX
Set-TextValue "E28" "  +8.29%  "

# Row 29 - PancakeSwap
Set-TextValue "E29" "  +2.07%  "

# Row 30 - RenderToken
Set-TextValue "E30" "  +4.77%  "

# Row 31 - EthereumClassic
Set-TextValue "D31" "25.76"
Set-TextValue "E31" "  +3.25%  "

# Row 32 - NEARProtocol
Set-TextValue "D32" "6.09"
Set-TextValue "E32" "  -0.47%  "

# Row 33 - Hedera
Set-TextValue "D33" "0.0984"
Set-TextValue "E33" "  -0.27%  "

# Row 34 - was Mantle, now PEPE
$ws.Range("B34").Value = "PEPE"
$ws.Range("C34").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue "D34" "0.0$([char]0x2083)0757"
Set-TextValue "E34" "  +18.55%  "

# Row 35 - was Filecoin, now Mantle
$ws.Range("B35").Value = "Mantle"
$ws.Range("C35").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue "D35" "0.973"
Set-TextValue "E35" "  +5.49%  "

# Row 36 - was PEPE, now Filecoin
$ws.Range("B36").Value = "Filecoin"
$ws.Range("C36").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D36" "5.77"
Set-TextValue "E36" "  +6.41%  "

# Row 37 - Stacks
Set-TextValue "D37" "2.10"
Set-TextValue "E37" "  +2.86%  "

# Row 38 - OKB
Set-TextValue "D38" "48.87"
Set-TextValue "E38" "  +1.36%  "

# Row 39 - Cosmos
Set-TextValue "D39" "8.85"
Set-TextValue "E39" "  +4.44%  "

# Row 40 - dogwifhat
Set-TextValue "D40" "2.75"
Set-TextValue "E40" "  +13.36%  "

# Row 41 - Bittensor
Set-TextValue "D41" "393.51"
Set-TextValue "E41" "  +7.51%  "

# Row 42 - was VeChain, now Kaspa
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D42" "0.108"
Set-TextValue "E42" "  +0.03%  "

# Row 43 - was Kaspa, now VeChain
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D43" "0.0349"
Set-TextValue "E43" "  -0.04%  "

# Row 44 - Maker
Set-TextValue "D44" "2.729.48"
Set-TextValue "E44" "  +3.56%  "

# Row 45 - TheGraph
Set-TextValue "D45" "0.247"
Set-TextValue "E45" "  +5.93%  "

# Row 46 - was USDe, now Monero
$ws.Range("B46").Value = "Monero"
$ws.Range("C46").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D46" "125.72"
Set-TextValue "E46" "  +5.68%  "

# Row 47 - was Monero, now USDe
$ws.Range("B47").Value = "USDe"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue "D47" "0.999"
Set-TextValue "E47" "  +0.01%  "

# Row 48 - Fetch.AI
Set-TextValue "D48" "2.03"
Set-TextValue "E48" "  +2.73%  "

# Row 49 - Stellar
Set-TextValue "E49" "  +2.12%  "

# Row 50 - InjectiveProtocol
Set-TextValue "D50" "23.37"
Set-TextValue "E50" "  +1.38%  "

# Row 51 - Arweave
Set-TextValue "D51" "31.95"
Set-TextValue "E51" "  +15.56%  "
